$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header in B1 changes from "Nombre humano" to "Sexo"
$ws.Range("B1").Value = "Sexo"

# Update the selected cell to E5 (as recorded in the saved file's sheetView)
$ws.Range("E5").Select()
